$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

# A1 is a brand-new cell (column A previously had no header); give it the
# same bold/centered/bordered "index" style already used by A2:A29 and the
# B1:M1 header row before writing its value.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$holeIds = @(
    "hole_id",
    "BRG_16_07",
    "BRG_16_03",
    "BRG_16_04B",
    "ECO_09_02",
    "BRG_16_01",
    "BRG_16_02",
    "BRG_01_05",
    "BRG_01_03",
    "BRG_05_09",
    "BRG_05_02",
    "ECO_09_01",
    "BRG_13_01",
    "BRG_16_09",
    "BRG_05_15",
    "BRG_05_13",
    "BRG_01_07",
    "BRG_05_05",
    "BRG_01_08",
    "ECO_09_05",
    "BRG_01_02",
    "BRG_05_11",
    "BRG_01_01",
    "BRG_05_10",
    "BRG_01_06",
    "BRG_05_03",
    "BRG_05_01",
    "BRG_01_09",
    "BRG_05_04"
)

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
